# Update to do list
# Remove the two rows for "Monster.GetAttackMessage" and
# "Dragon/Witch/Shulker/Skeleton/Warden Get Attack Message" (rows 18 & 19),
# shifting the remaining rows up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18:G19").EntireRow.Delete()

# Leave selection on B35, matching where the user ended up after editing.
$ws.Range("B35").Select()
